$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 201: Prashant
$ws.Range("B201").Value = "Prashant"
$ws.Range("C201").Value = 1
$ws.Range("D201").Value = 19
$ws.Range("E201").Value = 0
$ws.Range("F201").Value = 19
$ws.Range("G201").Value = 61050
$ws.Range("H201").Value = 80
$ws.Range("I201").Value = 40
$ws.Range("K201").Value = 365

# Row 202: Andy
$ws.Range("B202").Value = "Andy"
$ws.Range("C202").Value = 2
$ws.Range("D202").Value = 18
$ws.Range("E202").Value = 0
$ws.Range("F202").Value = 18
$ws.Range("G202").Value = 69100
$ws.Range("H202").Value = 60
$ws.Range("I202").Value = 20
$ws.Range("K202").Value = 349

# Row 203: Matt
$ws.Range("B203").Value = "Matt"
$ws.Range("C203").Value = 3
$ws.Range("D203").Value = 18
$ws.Range("E203").Value = 0
$ws.Range("F203").Value = 18
$ws.Range("G203").Value = 56500
$ws.Range("H203").Value = 60
$ws.Range("I203").Value = 20
$ws.Range("K203").Value = 362

# Row 204: Pepe
$ws.Range("B204").Value = "Pepe"
$ws.Range("C204").Value = 4
$ws.Range("D204").Value = 17
$ws.Range("E204").Value = 0
$ws.Range("F204").Value = 17
$ws.Range("G204").Value = 58800
$ws.Range("H204").Value = 30
$ws.Range("I204").Value = -10
$ws.Range("K204").Value = 364

# Row 205: Maisy
$ws.Range("B205").Value = "Maisy"
$ws.Range("C205").Value = 5
$ws.Range("D205").Value = 14
$ws.Range("E205").Value = 0
$ws.Range("F205").Value = 14
$ws.Range("G205").Value = 47400
$ws.Range("H205").Value = 20
$ws.Range("I205").Value = -20
$ws.Range("K205").Value = 360

# Row 206: Richard
$ws.Range("B206").Value = "Richard"
$ws.Range("C206").Value = 6
$ws.Range("D206").Value = 13
$ws.Range("E206").Value = 0
$ws.Range("F206").Value = 13
$ws.Range("G206").Value = 43700
$ws.Range("H206").Value = 80
$ws.Range("I206").Value = 40
$ws.Range("K206").Value = 366

# Row 207: Jon
$ws.Range("B207").Value = "Jon"
$ws.Range("C207").Value = 7
$ws.Range("D207").Value = 7
$ws.Range("E207").Value = 0
$ws.Range("F207").Value = 7
$ws.Range("G207").Value = 27450
$ws.Range("H207").Value = 0
$ws.Range("I207").Value = -40
$ws.Range("K207").Value = 357

# Row 208: Mark
$ws.Range("B208").Value = "Mark"
$ws.Range("C208").Value = 8
$ws.Range("D208").Value = 7
$ws.Range("E208").Value = 0
$ws.Range("F208").Value = 7
$ws.Range("G208").Value = 26450
$ws.Range("H208").Value = 10
$ws.Range("I208").Value = -20
$ws.Range("K208").Value = 361

# Row 209: Anthony
$ws.Range("B209").Value = "Anthony"
$ws.Range("C209").Value = 9
$ws.Range("D209").Value = 7
$ws.Range("E209").Value = 0
$ws.Range("F209").Value = 7
$ws.Range("G209").Value = 21000
$ws.Range("H209").Value = 10
$ws.Range("I209").Value = -10
$ws.Range("K209").Value = 350

# Row 210: Alex
$ws.Range("B210").Value = "Alex"
$ws.Range("C210").Value = 10
$ws.Range("D210").Value = 0
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 7700
$ws.Range("H210").Value = 0
$ws.Range("I210").Value = -20
$ws.Range("K210").Value = 348
